$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay TEXT even if it looks numeric
# (mirrors typing an apostrophe-prefixed entry in Excel), then strip the
# resulting quote-prefix formatting so no stray cell style is left behind.
function Set-TextValue($range, $text) {
    if ($text -match "^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$") {
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

Set-TextValue $ws.Range("D2") '26.667.84'
$ws.Range("E2").Value = '  +0.07%  '
Set-TextValue $ws.Range("D3") '1.598.96'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  +0.14%  '
Set-TextValue $ws.Range("D5") '211.47'
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("E7").Value = '  +0.12%  '
Set-TextValue $ws.Range("D8") '0.0619'
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E9").Value = '  +0.30%  '
Set-TextValue $ws.Range("D10") '19.59'
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("E11").Value = '  +0.52%  '
Set-TextValue $ws.Range("D12") '1.822.64'
$ws.Range("E12").Value = '  +0.35%  '
Set-TextValue $ws.Range("D13") '1.617.79'
$ws.Range("E13").Value = '  +1.54%  '
Set-TextValue $ws.Range("D14") '4.03'
$ws.Range("E14").Value = '  +0.29%  '
Set-TextValue $ws.Range("D15") '0.524'
$ws.Range("E15").Value = '  +0.37%  '
Set-TextValue $ws.Range("D16") '64.90'
$ws.Range("E16").Value = '  +0.28%  '
Set-TextValue $ws.Range("D17") '26.651.19'
$ws.Range("E17").Value = '  +0.11%  '
Set-TextValue $ws.Range("D18") '0.0₃0734'
$ws.Range("E18").Value = '  +0.65%  '
$ws.Range("E19").Value = '  +0.17%  '
Set-TextValue $ws.Range("D20") '208.08'
$ws.Range("E20").Value = '  -0.52%  '
Set-TextValue $ws.Range("D21") '7.07'
$ws.Range("E21").Value = '  +5.83%  '
Set-TextValue $ws.Range("D22") '4.28'
$ws.Range("E22").Value = '  +1.11%  '
Set-TextValue $ws.Range("D23") '2.33'
$ws.Range("E23").Value = '  +1.52%  '
Set-TextValue $ws.Range("D24") '8.94'
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("E25").Value = '  -0.85%  '
$ws.Range("E26").Value = '  +0.12%  '
Set-TextValue $ws.Range("D27") '7.13'
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("E28").Value = '  +0.07%  '
Set-TextValue $ws.Range("D29") '15.31'
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  +1.92%  '
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("E32").Value = '  +0.60%  '
$ws.Range("E33").Value = '  +1.13%  '
Set-TextValue $ws.Range("D34") '1.280.46'
$ws.Range("E34").Value = '  -0.98%  '
Set-TextValue $ws.Range("D35") '0.624'
$ws.Range("E35").Value = '  -8.43%  '
$ws.Range("E36").Value = '  +0.70%  '
Set-TextValue $ws.Range("D37") '1.49'
$ws.Range("E37").Value = '  +0.78%  '
Set-TextValue $ws.Range("D38") '0.0171'
$ws.Range("E38").Value = '  -0.40%  '
Set-TextValue $ws.Range("D39") '1.07'
$ws.Range("E39").Value = '  +19.90%  '
Set-TextValue $ws.Range("D40") '0.837'
$ws.Range("E40").Value = '  +0.04%  '
Set-TextValue $ws.Range("D41") '5.51'
$ws.Range("E41").Value = '  +3.05%  '
$ws.Range("E42").Value = '  +0.57%  '
$ws.Range("E43").Value = '  -0.52%  '
Set-TextValue $ws.Range("D44") '64.00'
$ws.Range("E44").Value = '  +0.81%  '
Set-TextValue $ws.Range("D45") '1.735.08'
$ws.Range("E45").Value = '  +0.42%  '
Set-TextValue $ws.Range("D46") '90.27'
$ws.Range("E46").Value = '  +0.71%  '
Set-TextValue $ws.Range("D47") '1.60'
$ws.Range("E47").Value = '  -2.44%  '
Set-TextValue $ws.Range("D48") '0.102'
$ws.Range("E48").Value = '  +3.63%  '
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("E51").Value = '  -0.68%  '
